# Append one new data row (row 71) to the LCY Saving Accounts sheet.
# Columns: A=TC, B=Account_ID, C=Customer_ID, D=PD
# New record: TC=118518, Account_ID=1008784273, Customer_ID=17866764, PD=6005
#
# These are ID-like strings that must be stored as text (not auto-converted
# to numbers), matching every other row already in the sheet. We briefly
# mark the target cells as Text before writing, then restore formatting so
# the cells end up with the default (unformatted) style, same as the rest
# of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 71
$target = $ws.Range("A" + $row + ":D" + $row)

$target.NumberFormat = "@"

$ws.Range("A" + $row).Value = "118518"
$ws.Range("B" + $row).Value = "1008784273"
$ws.Range("C" + $row).Value = "17866764"
$ws.Range("D" + $row).Value = "6005"

$target.ClearFormats()
